# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the Overview/zh-cn/de-de
# sheets move from "Handed back: in sync with en-US" to "Ready for handoff",
# the associated timestamps are refreshed, and the now-shorter status column
# is narrowed to fit.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws1.Range("G2").Value = "2016-09-05 01:06:17"

# ---- zh-cn sheet -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("H2").Value = "2016-09-05 01:06:12"

# ---- de-de sheet -------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("H2").Value = "2016-09-05 01:06:17"

# ---- Narrow the "Status" columns to fit the shorter text ---------------
# Target stored column width is ~17.216 characters. Excel's ColumnWidth
# property only takes effect in whole-pixel increments, so we feed it a
# value (16.3333...) that is guaranteed to land on the closest achievable
# pixel boundary to the target width.
$newStatusWidth = 16.3333333333333
$ws1.Columns.Item(5).ColumnWidth = $newStatusWidth   # Overview!E (zh-cn status)
$ws1.Columns.Item(6).ColumnWidth = $newStatusWidth   # Overview!F (de-de status)
$ws2.Columns.Item(3).ColumnWidth = $newStatusWidth   # zh-cn!C (Status)
$ws3.Columns.Item(3).ColumnWidth = $newStatusWidth   # de-de!C (Status)
